# opcodemeta.xlsx edit
#
# Commit message: "added detection of end of routine, added ToString() to
# provide disassembly, added more specificity to opcodes to help static
# analysis"
#
# Concretely (as seen in the sheet1.xml / sharedStrings.xml diff):
#   1. A new column is inserted at H ("end of routine?" flag column),
#      shifting the previous H..L (computed helper columns) one to the
#      right, to I..M. Excel auto-adjusts every formula reference when the
#      column is inserted, which reproduces the I/J/K/L/M formula shuffle
#      seen in the diff exactly.
#   2. The new column H is populated with the literal text "return" for
#      every opcode row that ends a routine (ret, rtrue, rfalse, print_ret,
#      ret_popped, quit) - i.e. "detection of end of routine". All other
#      rows are left blank, matching column G's look (Insert() already
#      carries the format across from G).
#   3. Move the selection to B3 to match the saved view state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new column before H; existing H:L shift to I:M and every
#    formula that referenced them is automatically rewritten by Excel.
$ws.Columns("H").Insert()

# 2) Flag the routine-terminating opcodes in the freshly inserted column.
#    Row -> opcode:
#      36 ret, 40 rtrue, 41 rfalse, 43 print_ret, 47 ret_popped, 49 quit
$returnRows = @(36, 40, 41, 43, 47, 49)
foreach ($r in $returnRows) {
    $ws.Range("H$r").Value = "return"
}

# 3) Restore the recorded selection/view state.
$ws.Range("B3").Select()
